$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = 2
$ws.Range("B22").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 45133
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 100112022
$ws.Range("G22").Value = "Arveja Verde"
$ws.Range("H22").Value = "Perfection"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 560
$ws.Range("K22").Value = 23000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 24000
$ws.Range("N22").Value = "$/malla 25 kilos"
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 960
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
